$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.054.46"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3
$ws.Range("D3").Value = "2.756.21"
$ws.Range("E3").Value = "  +1.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "578.39"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("D6").Value = "158.29"
$ws.Range("E6").Value = "  +2.41%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("E8").Value = "  -0.13%  "

# Row 10
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.386"
$ws.Range("E10").Value = "  -1.08%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  -15.02%  "

# Row 12
$ws.Range("E12").Value = "  -2.22%  "

# Row 13
$ws.Range("D13").Value = "3.244.52"
$ws.Range("E13").Value = "  +1.01%  "

# Row 14
$ws.Range("D14").Value = "26.88"
$ws.Range("E14").Value = "  +2.16%  "

# Row 15
$ws.Range("D15").Value = "63.760.50"
$ws.Range("E15").Value = "  -0.08%  "

# Row 16
$ws.Range("E16").Value = "  -0.29%  "

# Row 17
$ws.Range("D17").Value = "2.760.07"
$ws.Range("E17").Value = "  +0.70%  "

# Row 18
$ws.Range("D18").Value = "12.17"
$ws.Range("E18").Value = "  +1.86%  "

# Row 19
$ws.Range("D19").Value = "4.87"
$ws.Range("E19").Value = "  +0.22%  "

# Row 20
$ws.Range("D20").Value = "359.51"
$ws.Range("E20").Value = "  -0.31%  "

# Row 21
$ws.Range("E21").Value = "  -1.60%  "

# Row 22
$ws.Range("E22").Value = "  +2.54%  "

# Row 23
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").Value = "65.76"
$ws.Range("E24").Value = "  -0.49%  "

# Row 25
$ws.Range("E25").Value = "  +1.21%  "

# Row 26
$ws.Range("D26").Value = "8.51"
$ws.Range("E26").Value = "  -0.27%  "

# Row 27
$ws.Range("E27").Value = "  +0.03%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("E28").Value = "  +1.90%  "

# Row 29
$ws.Range("E29").Value = "  -1.48%  "

# Row 30
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").Value = "  -1.17%  "

# Row 31
$ws.Range("E31").Value = "  +0.45%  "

# Row 32
$ws.Range("D32").Value = "167.79"
$ws.Range("E32").Value = "  -2.19%  "

# Row 33
$ws.Range("D33").Value = "20.32"
$ws.Range("E33").Value = "  -0.77%  "

# Row 34
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  +3.37%  "

# Row 35
$ws.Range("E35").Value = "  +0.10%  "

# Row 37
$ws.Range("E37").Value = "  -0.40%  "

# Row 38
$ws.Range("D38").Value = "0.992"
$ws.Range("E38").Value = "  -0.42%  "

# Row 39
$ws.Range("D39").Value = "6.25"
$ws.Range("E39").Value = "  +12.28%  "

# Row 40
$ws.Range("D40").Value = "4.16"
$ws.Range("E40").Value = "  -1.25%  "

# Row 41
$ws.Range("D41").Value = "329.55"

# Row 42
$ws.Range("D42").Value = "39.34"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("D43").Value = "21.57"
$ws.Range("E43").Value = "  -0.52%  "

# Row 44
$ws.Range("D44").Value = "0.0594"
$ws.Range("E44").Value = "  +0.62%  "

# Row 45
$ws.Range("D45").Value = "21.79"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("E46").Value = "  +1.04%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "136.80"
$ws.Range("E47").Value = "  -1.95%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.634"
$ws.Range("E48").Value = "  -1.36%  "

# Row 49
$ws.Range("E49").Value = "  +0.71%  "

# Row 50
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.14%  "

# Row 51
$ws.Range("E51").Value = "  +0.63%  "
